$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Append a new row of data (2021/09/15 report) right after the last existing row (80)
$row = 81
$prev = $row - 1

# Remember the number formats used by the row above so the new row matches it
$fmtA = $ws.Cells.Item($prev, 1).NumberFormat
$fmtNum = $ws.Cells.Item($prev, 2).NumberFormat

# Column A holds date-like labels stored as plain text (not real dates), so force
# a Text format first - otherwise Excel auto-converts the "2021/09/15" string into
# a date serial number. Revert the format back afterwards so the cell keeps the
# same style as the rest of the column.
$ws.Cells.Item($row, 1).NumberFormat = "@"
$ws.Cells.Item($row, 1).Value = "2021/09/15"
$ws.Cells.Item($row, 1).NumberFormat = $fmtA

$ws.Cells.Item($row, 2).NumberFormat = $fmtNum
$ws.Cells.Item($row, 2).Value = 191.1
$ws.Cells.Item($row, 3).NumberFormat = $fmtNum
$ws.Cells.Item($row, 3).Value = 196.1
$ws.Cells.Item($row, 4).NumberFormat = $fmtNum
$ws.Cells.Item($row, 4).Value = 0.84
$ws.Cells.Item($row, 5).NumberFormat = $fmtNum
$ws.Cells.Item($row, 5).Value = 0.83

# Update the active selection to the cell right after the freshly entered row,
# matching how Excel leaves the selection after typing in a new row of data.
$ws.Range("A" + ($row + 1)).Select()
